$wb = $excel.ActiveWorkbook

# --- "finance" sheet ---
$ws = $wb.Worksheets.Item("finance")

# Rename the API label cells to start with "*", same treatment as the
# existing "*Note" / "*Args" labels already used on this sheet.
$ws.Range("A2").Value = "*Api"
$ws.Range("A4").Value = "*Args"

# Duplicate the existing "Api query / Args table" block (rows 2-5) into a
# second block in rows 6-9 describing the route-generated "query1" api.
$ws.Range("A6").Value = "*Api"
$ws.Range("B6").Value = "query1"

$ws.Range("A7").Value = "*Note"
$ws.Range("B7").Value = "args_name"
$ws.Range("C7").Value = "type"
$ws.Range("D7").Value = "required"
$ws.Range("E7").Value = "missing"

$ws.Range("A8").Value = "*Args"
$ws.Range("B8").Value = "ops_org_id"
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = "xxx"

$ws.Range("B9").Value = "page_no"
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 10

# The blank cells around/below the tables pick up the same "0" number
# format already used by the populated numeric columns.
$ws.Range("D1:F2").NumberFormat = "0"
$ws.Range("F3:F5").NumberFormat = "0"
$ws.Range("A5").NumberFormat = "0"
$ws.Range("A6:F11").NumberFormat = "0"

# --- "_Note" sheet ---
$ws2 = $wb.Worksheets.Item("_Note")
$ws2.Range("B1:E1").NumberFormat = "0"
$ws2.Range("A2:E10").NumberFormat = "0"
